$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$name = $ws.Range("A2").Text
$cmc = $ws.Range("A3").Text
$type = $ws.Range("A4").Text
$ability1 = $ws.Range("A5").Text
$ability2 = $ws.Range("A6").Text
$ability3 = $ws.Range("A7").Text
$pt = $ws.Range("A8").Text

$combined = "('" + $name + "', ['" + $cmc + "', '" + $type + "', '" + $ability1 + "', '" + $ability2 + "', '" + $ability3 + "', '" + $pt + "'])"

$ws.Range("A2").Value = $combined

$ws.Range("A3:A8").EntireRow.Delete()
